# Weekly fruit/vegetable price update.
# Inserts two new price records (row 100 and 101) into the daily logic
# sub-workbook for Ají (Vega Central Mapocho de Santiago), pushing the
# existing rows 100-157 down to 102-159.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 100, shifting all
# subsequent rows (old 100-157) down to 102-159.
$ws.Rows.Item(100).Resize(2).Insert()

# New row 100: Ají / Inferno / Primera, Región de Arica y Parinacota
$ws.Cells.Item(100, 1).Value = 9
$ws.Cells.Item(100, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(100, 3).Value = "Metropolitana"
$ws.Cells.Item(100, 4).Value = 44452
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = 100112021
$ws.Cells.Item(100, 7).Value = "Ají"
$ws.Cells.Item(100, 8).Value = "Inferno"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 25
$ws.Cells.Item(100, 11).Value = 48000
$ws.Cells.Item(100, 12).Value = 50000
$ws.Cells.Item(100, 13).Value = 48960
$ws.Cells.Item(100, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(100, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(100, 16).Value = 4080
$ws.Cells.Item(100, 17).Value = 12
$ws.Cells.Item(100, 18).Value = "Hortaliza"

# New row 101: Ají / Inferno / Segunda, Región de Arica y Parinacota
$ws.Cells.Item(101, 1).Value = 9
$ws.Cells.Item(101, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(101, 3).Value = "Metropolitana"
$ws.Cells.Item(101, 4).Value = 44452
$ws.Cells.Item(101, 5).Value = 13
$ws.Cells.Item(101, 6).Value = 100112021
$ws.Cells.Item(101, 7).Value = "Ají"
$ws.Cells.Item(101, 8).Value = "Inferno"
$ws.Cells.Item(101, 9).Value = "Segunda"
$ws.Cells.Item(101, 10).Value = 7
$ws.Cells.Item(101, 11).Value = 45000
$ws.Cells.Item(101, 12).Value = 47000
$ws.Cells.Item(101, 13).Value = 46143
$ws.Cells.Item(101, 14).Value = "`$/caja 12 kilos"
$ws.Cells.Item(101, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(101, 16).Value = 3845
$ws.Cells.Item(101, 17).Value = 12
$ws.Cells.Item(101, 18).Value = "Hortaliza"

Write-Host "Inserted rows 100-101; sheet now spans $($ws.UsedRange.Address())"
